$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "a0".."a3" pin-name labels in column C (rows 10-13) to "d10".."d13"
$ws.Range("C10").Value = "d10"
$ws.Range("C11").Value = "d11"
$ws.Range("C12").Value = "d12"
$ws.Range("C13").Value = "d13"

# Move the active selection from D11 to C14
$ws.Range("C14").Select()
